# إضافة حدث جديد في Card21 by admin at 2025-12-08 08:36:18
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card21")
$ws.Activate()

# Row 17: columns B..K were blank placeholders; author filled them with the
# literal "nan" text (matches the rest of the sheet's "not applicable" marker).
$ws.Range("B17:K17").Value = "nan"

# Row 18: brand-new service-log entry for Card21.
$ws.Range("A18").Value = "21"
$ws.Range("B18:K18").Value = ""
$ws.Range("L18").Value = "21\1\2025"
$ws.Range("M18").Value = ""
$ws.Range("N18").Value = "تم سن الفلاتس والسليندر وتغيير الجرائد الخلفيه ومعياره المكنه (1_5_8)"
$ws.Range("O18").Value = "الخبير"
